$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mal")
$ws.Rows.Item(3).Insert()
$ws.Range("C3").Value = "MEMOday3"
$ws.Range("D3").Value = "313805"
$ws.Range("E7").Select()
